$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NATI")

# Insert a new column before D; this shifts existing D:K data to E:L
# and keeps each row's cell style uniform because Excel carries the
# left-neighbor formatting into the new column by default. We then
# explicitly copy number-format/style from column E (the former D)
# into the new D for the rows that actually hold financial data, and
# finally set the copied-forward width.
$ws.Columns("D:D").Insert()

# Re-apply D-column formatting (number format + font) from column E
# for each contiguous data block (header rows 5,6,37,79 have no D cell).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep column D the same width as its neighbors post-insert.
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Populate the new column D with the newest fiscal-year figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 1359100
$ws.Range("D9").Value2 = 333700
$ws.Range("D10").Value2 = 1025400
$ws.Range("D12").Value2 = 258900
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 13900
$ws.Range("D15").Value2 = 2200
$ws.Range("D17").Value2 = 1185500
$ws.Range("D18").Value2 = 173600
$ws.Range("D20").Value2 = 2900
$ws.Range("D21").Value2 = 247100
$ws.Range("D22").Value2 = 0
$ws.Range("D23").Value2 = 176500
$ws.Range("D24").Value2 = 25600
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 150900
$ws.Range("D27").Value2 = 150900
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 4200
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -2900
$ws.Range("D33").Value2 = 155100
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 155100
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 259400
$ws.Range("D42").Value2 = 271400
$ws.Range("D43").Value2 = 243000
$ws.Range("D44").Value2 = 194100
$ws.Range("D45").Value2 = 54300
$ws.Range("D46").Value2 = 1022200
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 245200
$ws.Range("D49").Value2 = 375300
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 28500
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 1671200
$ws.Range("D57").Value2 = 48400
$ws.Range("D58").Value2 = 0
$ws.Range("D59").Value2 = 234600
$ws.Range("D60").Value2 = 283000
$ws.Range("D61").Value2 = 0
$ws.Range("D62").Value2 = 149900
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 432900
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 356400
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 1238400
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 155100
$ws.Range("D83").Value2 = 70700
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 274600
$ws.Range("D91").Value2 = -34700
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -210000
$ws.Range("D96").Value2 = -121500
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -90800
$ws.Range("D101").Value2 = -4500
$ws.Range("D102").Value2 = -30800

# Row 91 (Capital Expenditures) received corrected historical figures,
# not just a column shift, for years 2012-2018 (columns D:J); column K
# already holds the correct shifted value (-71900) after the insert.
$ws.Range("D91").Value2 = -34700
$ws.Range("E91").Value2 = -30300
$ws.Range("F91").Value2 = -44400
$ws.Range("G91").Value2 = -34000
$ws.Range("H91").Value2 = -44900
$ws.Range("I91").Value2 = -47800
$ws.Range("J91").Value2 = -89100
